{"js": "// Replace each three-digit-divided-by-one-digit problem in the worksheet\n// table with the new problem from the updated master. Every \"before\"\n// string is unique within the document, so an exact, case-sensitive\n// search safely targets the single matching run regardless of order.\nconst replacements = [\n  [\"947\u00f72=\", \"246\u00f79=\"],\n  [\"768\u00f78=\", \"165\u00f72=\"],\n  [\"331\u00f77=\", \"523\u00f73=\"],\n  [\"927\u00f77=\", \"321\u00f74=\"],\n  [\"352\u00f75=\", \"134\u00f72=\"],\n  [\"354\u00f75=\", \"582\u00f72=\"],\n  [\"557\u00f78=\", \"885\u00f76=\"],\n  [\"581\u00f75=\", \"596\u00f79=\"],\n  [\"350\u00f74=\", \"331\u00f77=\"],\n  [\"614\u00f73=\", \"783\u00f77=\"],\n  [\"589\u00f79=\", \"214\u00f77=\"],\n  [\"162\u00f76=\", \"841\u00f76=\"],\n  [\"104\u00f79=\", \"690\u00f79=\"],\n  [\"632\u00f77=\", \"227\u00f78=\"],\n  [\"546\u00f75=\", \"733\u00f72=\"],\n  [\"760\u00f76=\", \"696\u00f73=\"],\n  [\"647\u00f78=\", \"505\u00f79=\"],\n  [\"855\u00f78=\", \"263\u00f72=\"],\n  [\"369\u00f77=\", \"829\u00f76=\"],\n  [\"137\u00f76=\", \"790\u00f72=\"],\n  [\"781\u00f72=\", \"559\u00f75=\"],\n  [\"970\u00f72=\", \"253\u00f78=\"],\n  [\"344\u00f77=\", \"576\u00f76=\"],\n  [\"379\u00f77=\", \"408\u00f74=\"],\n  [\"800\u00f79=\", \"245\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-divided-by-one-digit problem in the worksheet\n# table with the new problem from the updated master. Every \"before\"\n# string is unique within the document, so Find/Replace targets exactly\n# one run per pair regardless of the order the pairs are processed in.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"947\u00f72=\", \"246\u00f79=\"),\n  @(\"768\u00f78=\", \"165\u00f72=\"),\n  @(\"331\u00f77=\", \"523\u00f73=\"),\n  @(\"927\u00f77=\", \"321\u00f74=\"),\n  @(\"352\u00f75=\", \"134\u00f72=\"),\n  @(\"354\u00f75=\", \"582\u00f72=\"),\n  @(\"557\u00f78=\", \"885\u00f76=\"),\n  @(\"581\u00f75=\", \"596\u00f79=\"),\n  @(\"350\u00f74=\", \"331\u00f77=\"),\n  @(\"614\u00f73=\", \"783\u00f77=\"),\n  @(\"589\u00f79=\", \"214\u00f77=\"),\n  @(\"162\u00f76=\", \"841\u00f76=\"),\n  @(\"104\u00f79=\", \"690\u00f79=\"),\n  @(\"632\u00f77=\", \"227\u00f78=\"),\n  @(\"546\u00f75=\", \"733\u00f72=\"),\n  @(\"760\u00f76=\", \"696\u00f73=\"),\n  @(\"647\u00f78=\", \"505\u00f79=\"),\n  @(\"855\u00f78=\", \"263\u00f72=\"),\n  @(\"369\u00f77=\", \"829\u00f76=\"),\n  @(\"137\u00f76=\", \"790\u00f72=\"),\n  @(\"781\u00f72=\", \"559\u00f75=\"),\n  @(\"970\u00f72=\", \"253\u00f78=\"),\n  @(\"344\u00f77=\", \"576\u00f76=\"),\n  @(\"379\u00f77=\", \"408\u00f74=\"),\n  @(\"800\u00f79=\", \"245\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
